$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing header-ish cells (question/answer -> question1/amswer1)
$ws.Range("A1").Value = "question1"
$ws.Range("B1").Value = "amswer1"

# Add new rows 2 and 3 with question2/amswer2 pairs
$ws.Range("A2").Value = "question2"
$ws.Range("B2").Value = "amswer2"
$ws.Range("A3").Value = "question2"
$ws.Range("B3").Value = "amswer2"
